$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.178.79"
$ws.Range("E2").Value = "  +0.54%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.664.05"
$ws.Range("E3").Value = "  +3.47%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D5").Value = "'607.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.57%  "

# Row 6 - Solana (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D6").Value = "'143.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.67%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.15%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.664.44"
$ws.Range("E9").Value = "  +3.51%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.01%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +1.53%  "

# Row 12 - TRON (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D12").Value = "'0.152"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.91%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  +1.95%  "

# Row 14 - Avalanche
$ws.Range("E14").Value = "  +1.11%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.142.51"
$ws.Range("E15").Value = "  +3.55%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "63.057.60"
$ws.Range("E16").Value = "  +0.50%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +0.01%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.681.78"
$ws.Range("E18").Value = "  +4.46%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  +3.25%  "

# Row 20 - BitcoinCash (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D20").Value = "'338.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.53%  "

# Row 21 - Polkadot
$ws.Range("E21").Value = "  +1.26%  "

# Row 22 - Uniswap (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D22").Value = "'6.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.34%  "

# Row 23 - Dai (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "

# Row 24 - Litecoin (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D24").Value = "'67.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  +3.34%  "

# Row 26 - SuiNetwork
$ws.Range("E26").Value = "  -2.44%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  +0.17%  "

# Row 28 - now Binance-PegBSC-USD (was InternetComputer(DFINITY))
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.45%  "

# Row 29 - now InternetComputer(DFINITY) (was Binance-PegBSC-USD)
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'8.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.95%  "

# Row 30 - Bittensor (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D30").Value = "'536.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +18.31%  "

# Row 31 - Aptos
$ws.Range("E31").Value = "  -1.61%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +5.28%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  +9.51%  "

# Row 34 - PEPE
$ws.Range("D34").Value = "0.0₃0809"
$ws.Range("E34").Value = "  +1.50%  "

# Row 35 - Monero (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D35").Value = "'173.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.52%  "

# Row 36 - NEARProtocol (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D36").Value = "'5.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +14.28%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  +0.04%  "

# Row 38 - PolygonEcosystemToken
$ws.Range("E38").Value = "  +1.44%  "

# Row 39 - EthereumClassic (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D39").Value = "'19.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.75%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  +7.59%  "

# Row 41 - Aave (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D41").Value = "'174.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.50%  "

# Row 43 - Filecoin (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D43").Value = "'3.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "

# Row 44 - InjectiveProtocol
$ws.Range("E44").Value = "  +5.12%  "

# Row 45 - Hedera
$ws.Range("E45").Value = "  +5.00%  "

# Row 46 - Mantle (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D46").Value = "'0.632"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.24%  "

# Row 47 - Stellar
$ws.Range("E47").Value = "  +0.07%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  +1.74%  "

# Row 49 - EnergySwap (D value is otherwise-numeric, force text with leading apostrophe)
$ws.Range("D49").Value = "'18.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.14%  "

# Row 51 - WhiteBITCoin
$ws.Range("E51").Value = "  -0.72%  "
